# Fix spelling error: "extrernaldoor" -> "externaldoor" in the Device Code column (E10:E12)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "externaldoor"
$ws.Range("E11").Value = "externaldoor"
$ws.Range("E12").Value = "externaldoor"

# Update the selection to reflect the cells that were just corrected
$ws.Range("E10:E12").Select()
